# Adding open extent report feature
# - Row 3 (the "SecondScript"/"DS002" test row) now also reports that the
#   HTML extent report link is available, matching row 2's "Yes" in column C.
# - Selection moves to A3 to reflect the newly active/edited row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Yes"
$ws.Range("A3").Select()
